$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Ají" (Hortaliza) at
# Terminal La Palmera de La Serena. Insert a new row at position 265 so
# that existing rows 265-323 shift down to 266-324 (same as Excel's
# "Insert Copied/Shifted Cells" behaviour when a row is inserted above).
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A265").Value = 8
$ws.Range("B265").Value = "Terminal La Palmera de La Serena"
$ws.Range("C265").Value = "Coquimbo"
$ws.Range("D265").Value = 44889
$ws.Range("E265").Value = 4
$ws.Range("F265").Value = 100112021
$ws.Range("G265").Value = "Ají"
$ws.Range("H265").Value = "Inferno"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 400
$ws.Range("K265").Value = 29000
$ws.Range("L265").Value = 30000
$ws.Range("M265").Value = 29500
$ws.Range("N265").Value = "`$/caja 25 kilos"
$ws.Range("O265").Value = "Provincia de Limarí"
$ws.Range("P265").Value = 1180
$ws.Range("Q265").Value = 25
$ws.Range("R265").Value = "Hortaliza"
